# Update "想去人数" (want-to-go count) figures to the freshly scraped values.
# Each (sheet, cell) pair below corresponds to one event whose counter
# increased by 1 since the last site build.
$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Sheet = "展览";     Cell = "F28"; Value = 36 },
    @{ Sheet = "展览";     Cell = "F31"; Value = 327 },
    @{ Sheet = "展览";     Cell = "F38"; Value = 4030 },
    @{ Sheet = "演出";     Cell = "F6";  Value = 195 },
    @{ Sheet = "演出";     Cell = "F22"; Value = 265 },
    @{ Sheet = "本地生活"; Cell = "F4";  Value = 1286 },
    @{ Sheet = "全部类型"; Cell = "F2";  Value = 1286 },
    @{ Sheet = "全部类型"; Cell = "F15"; Value = 195 },
    @{ Sheet = "全部类型"; Cell = "F16"; Value = 195 },
    @{ Sheet = "全部类型"; Cell = "F32"; Value = 36 },
    @{ Sheet = "全部类型"; Cell = "F34"; Value = 265 },
    @{ Sheet = "全部类型"; Cell = "F37"; Value = 327 },
    @{ Sheet = "全部类型"; Cell = "F49"; Value = 4030 }
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    $ws.Range($u.Cell).Value = $u.Value
}
